# Powerpoint writer: consolidate text runs when possible.
# Rewrites each multi-run paragraph's text as a single run by clearing the
# existing text range and re-inserting the combined string.

function Set-ConsolidatedText($shape, [string]$text) {
    $tr = $shape.TextFrame.TextRange
    [void]$tr.Delete()
    $tr.InsertAfter($text) | Out-Null
}

$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
Set-ConsolidatedText $s1.Shapes.Item(1) "Slide 1"
Set-ConsolidatedText $s1.Shapes.Item(3) "an image"

$s2 = $p.Slides.Item(2)
Set-ConsolidatedText $s2.Shapes.Item(1) "Slide 2"
Set-ConsolidatedText $s2.Shapes.Item(4) "an image"
